$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the quarter-final results (columns J = Home_Score, K = Away_Score) ---
$ws.Range("J46").Value = 1
$ws.Range("K46").Value = 1

$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0

$ws.Range("J48").Value = 1
$ws.Range("K48").Value = 1

$ws.Range("J49").Value = 2
$ws.Range("K49").Value = 1

# --- Add the two new semi-final fixture rows ---
$ws.Range("A50").Value = "Tue"
$ws.Range("B50").Value = "Jul 09, 2024"
$ws.Range("C50").Value = "21:00:00"
$ws.Range("D50").Value = "Spain"
$ws.Range("G50").Value = "France"
$ws.Range("H50").Value = "Munich"

$ws.Range("A51").Value = "Wed"
$ws.Range("B51").Value = "Jul 10, 2024"
$ws.Range("C51").Value = "21:00:00"
$ws.Range("D51").Value = "Netherlands"
$ws.Range("G51").Value = "England"
$ws.Range("H51").Value = "Dortmund"

# --- Widen column B (Date) slightly, as in the authored workbook ---
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666

# --- Restore the scroll position / selection the author left the sheet in ---
$excel.ActiveWindow.ScrollRow = 36
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C58").Select()
